# Automatic update of files.
# Increment the "Taxonsorteringsordning" (column B) value by 1 for each
# data row (rows 2-11) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $cell.Value2 + 1
}
